# Adds 14 new break-log entries (rows 141-154) for 2026-01-27 that were
# captured after the workbook's last export, and extends the sheet's
# used range accordingly (A1:H140 -> A1:H154).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=141; A=7025939006; B='Booster_lav'; C='🎀 𝒽𝒾_𝓁𝒶𝒶𝓋 🎀'; D='🚻 Comfort Room'; E='OUT'; F='2026-01-27 10:27:29'; G=$null; H=$null },
    @{ Row=142; A=7025939006; B='Booster_lav'; C='🎀 𝒽𝒾_𝓁𝒶𝒶𝓋 🎀'; D='🚻 Comfort Room'; E='BACK'; F='2026-01-27 10:31:11'; G=3.7; H=$null },
    @{ Row=143; A=8313813326; B='booster_roxan'; C='roxy'; D='🚻 Comfort Room'; E='BACK'; F='2026-01-27 10:33:02'; G=185.5; H=$null },
    @{ Row=144; A=8397936341; B='Booster_Moja'; C='ʍօʝǟ🔫🔥'; D='🚬 Smoke Break'; E='OUT'; F='2026-01-27 11:00:02'; G=$null; H=$null },
    @{ Row=145; A=8203583816; B='Cyrus0228'; C='Cyrus Rufo'; D='🚬 Smoke Break'; E='OUT'; F='2026-01-27 11:02:44'; G=$null; H=$null },
    @{ Row=146; A=8224136102; B='Matiluk'; C='Boost-Lo'; D='🚬 Smoke Break'; E='OUT'; F='2026-01-27 14:31:20'; G=$null; H=$null },
    @{ Row=147; A=8224136102; B='Matiluk'; C='Boost-Lo'; D='🚬 Smoke Break'; E='BACK'; F='2026-01-27 14:31:24'; G=0.1; H=$null },
    @{ Row=148; A=8224136102; B='Matiluk'; C='Boost-Lo'; D='🚻 Comfort Room'; E='OUT'; F='2026-01-27 15:36:06'; G=$null; H=$null },
    @{ Row=149; A=8224136102; B='Matiluk'; C='Boost-Lo'; D='🚻 Comfort Room'; E='BACK'; F='2026-01-27 16:41:26'; G=65.3; H=$null },
    @{ Row=150; A=8013843575; B='Booster_yham'; C='yham'; D='🚬 Smoke Break'; E='OUT'; F='2026-01-27 18:37:50'; G=$null; H=$null },
    @{ Row=151; A=8013843575; B='Booster_yham'; C='yham'; D='🚬 Smoke Break'; E='BACK'; F='2026-01-27 18:55:54'; G=18.1; H=$null },
    @{ Row=152; A=8011222190; B='syintel'; C='Sheena'; D='🚻 Comfort Room'; E='OUT'; F='2026-01-27 20:39:06'; G=$null; H=$null },
    @{ Row=153; A=8011222190; B='syintel'; C='Sheena'; D='🚻 Comfort Room'; E='BACK'; F='2026-01-27 20:42:17'; G=3.2; H=$null },
    @{ Row=154; A=8011222190; B='syintel'; C='Sheena'; D='🚻 Comfort Room'; E='OUT'; F='2026-01-27 21:27:40'; G=$null; H=$null }
)

foreach ($row in $rows) {
    $r = $row.Row
    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F

    if ($row.G -ne $null) {
        $ws.Range("G$r").Value = $row.G
    } else {
        $ws.Range("G$r").Value = ""
    }

    if ($row.H -ne $null) {
        $ws.Range("H$r").Value = $row.H
    } else {
        $ws.Range("H$r").Value = ""
    }
}
